$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '30.200.52'
Set-TextValue $ws.Range('E2') '  +0.13%  '
Set-TextValue $ws.Range('D3') '1.858.41'
Set-TextValue $ws.Range('D4') '1.000'
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '235.90'
Set-TextValue $ws.Range('E5') '  +0.79%  '
Set-TextValue $ws.Range('D6') '1.000'
Set-TextValue $ws.Range('E6') '  -0.01%  '
Set-TextValue $ws.Range('D7') '0.4669'
Set-TextValue $ws.Range('E7') '  +0.02%  '
Set-TextValue $ws.Range('E8') '  +0.92%  '
Set-TextValue $ws.Range('D9') '0.06513'
Set-TextValue $ws.Range('E9') '  -0.72%  '
Set-TextValue $ws.Range('D10') '21.66'
Set-TextValue $ws.Range('E10') '  +8.36%  '
Set-TextValue $ws.Range('D11') '0.07899'
Set-TextValue $ws.Range('E11') '  +0.90%  '
Set-TextValue $ws.Range('D12') '97.22'
Set-TextValue $ws.Range('E12') '  +0.39%  '
Set-TextValue $ws.Range('D13') '1.865.09'
Set-TextValue $ws.Range('E13') '  +0.38%  '
Set-TextValue $ws.Range('D14') '5.150'
Set-TextValue $ws.Range('E14') '  +0.76%  '
Set-TextValue $ws.Range('D15') '0.6775'
Set-TextValue $ws.Range('E15') '  +1.91%  '
Set-TextValue $ws.Range('D16') '278.76'
Set-TextValue $ws.Range('E16') '  -1.22%  '
Set-TextValue $ws.Range('D17') '30.201.63'
Set-TextValue $ws.Range('E17') '  -0.26%  '
Set-TextValue $ws.Range('D18') '13.56'
Set-TextValue $ws.Range('E18') '  +7.68%  '
Set-TextValue $ws.Range('E19') '  +0.01%  '
Set-TextValue $ws.Range('D20') '5.368'
Set-TextValue $ws.Range('E20') '  -1.22%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D21') '2.109.44'
Set-TextValue $ws.Range('E21') '  -0.09%  '
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D22') '0.000007300'
Set-TextValue $ws.Range('E22') '  +0.89%  '
Set-TextValue $ws.Range('D23') '1.000'
Set-TextValue $ws.Range('E23') '  -0.05%  '
Set-TextValue $ws.Range('D24') '6.142'
Set-TextValue $ws.Range('E24') '  +0.03%  '
Set-TextValue $ws.Range('D25') '166.59'
Set-TextValue $ws.Range('E25') '  -0.51%  '
Set-TextValue $ws.Range('D26') '9.179'
Set-TextValue $ws.Range('D27') '19.04'
Set-TextValue $ws.Range('E27') '  +0.67%  '
Set-TextValue $ws.Range('D28') '1.921'
Set-TextValue $ws.Range('E28') '  +0.20%  '
Set-TextValue $ws.Range('D29') '1.382'
Set-TextValue $ws.Range('E29') '  +3.17%  '
Set-TextValue $ws.Range('D30') '0.09699'
Set-TextValue $ws.Range('E30') '  +1.42%  '
Set-TextValue $ws.Range('D31') '4.362'
Set-TextValue $ws.Range('E31') '  -0.99%  '
Set-TextValue $ws.Range('D32') '1.474'
Set-TextValue $ws.Range('E32') '  +0.30%  '
Set-TextValue $ws.Range('D33') '4.029'
Set-TextValue $ws.Range('E33') '  -1.65%  '
Set-TextValue $ws.Range('D34') '0.04706'
Set-TextValue $ws.Range('E34') '  +1.09%  '
Set-TextValue $ws.Range('D35') '1.129'
Set-TextValue $ws.Range('E35') '  +2.83%  '
Set-TextValue $ws.Range('D36') '0.7045'
Set-TextValue $ws.Range('E36') '  +0.66%  '
Set-TextValue $ws.Range('D37') '2.707'
Set-TextValue $ws.Range('E37') '  -0.06%  '
Set-TextValue $ws.Range('D38') '0.01857'
Set-TextValue $ws.Range('E38') '  +0.36%  '
Set-TextValue $ws.Range('D39') '2.625'
Set-TextValue $ws.Range('E39') '  +4.61%  '
Set-TextValue $ws.Range('D40') '6.305'
Set-TextValue $ws.Range('E40') '  -0.79%  '
Set-TextValue $ws.Range('D41') '74.32'
Set-TextValue $ws.Range('E41') '  +3.27%  '
Set-TextValue $ws.Range('D42') '1.944'
Set-TextValue $ws.Range('E42') '  +0.72%  '
Set-TextValue $ws.Range('D43') '0.8470'
Set-TextValue $ws.Range('E43') '  -0.72%  '
Set-TextValue $ws.Range('D44') '0.9999'
Set-TextValue $ws.Range('E44') '  -0.06%  '
Set-TextValue $ws.Range('D45') '0.4160'
Set-TextValue $ws.Range('E45') '  +0.06%  '
Set-TextValue $ws.Range('E46') '  -0.34%  '
Set-TextValue $ws.Range('D47') '984.01'
Set-TextValue $ws.Range('E47') '  -0.46%  '
Set-TextValue $ws.Range('D48') '7.145'
Set-TextValue $ws.Range('E48') '  -0.95%  '
Set-TextValue $ws.Range('D49') '9.252'
Set-TextValue $ws.Range('E49') '  +0.79%  '
Set-TextValue $ws.Range('D50') '33.98'
Set-TextValue $ws.Range('E50') '  +0.14%  '
Set-TextValue $ws.Range('D51') '0.05642'
Set-TextValue $ws.Range('E51') '  +0.15%  '
